# demo-02-advanced.xlsx : add a "Year" header column to both tables on
# Demo2 (left table B4:D13, right table F4:H13), matching the bold,
# bordered, centered header style already used for Name/Date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demo2")

# --- New "Year" header label for both tables --------------------------
$ws.Range("B4").Value = "Year"
$ws.Range("F4").Value = "Year"

# --- Header row (B4:D4 and F4:H4) gets bold text, a thin box border and
#     centered text, same look for the pre-existing Name/Date cells too.
$leftHeader = $ws.Range("B4:D4")
$rightHeader = $ws.Range("F4:H4")

$leftHeader.Font.Bold = $true
$rightHeader.Font.Bold = $true

$leftHeader.Borders.LineStyle = 1
$leftHeader.Borders.Weight = 2
$rightHeader.Borders.LineStyle = 1
$rightHeader.Borders.Weight = 2

$leftHeader.HorizontalAlignment = -4108
$rightHeader.HorizontalAlignment = -4108

# --- Cursor/selection moved one column left (J3 -> I3) -----------------
$ws.Range("I3").Select()
